$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift product-row content (name / balance / price / transactions) from
#     rows 25..30 down to 26..31, bottom-up so we never clobber data we still need.
#     NOTE: column A is just the running sequence number tied to the row position
#     (row - 3) and does not move with the product data. ---
for ($r = 30; $r -ge 25; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 2).Value()  = $ws.Cells.Item($r, 2).Value()   # B - name
    $ws.Cells.Item($dst, 8).Value()  = $ws.Cells.Item($r, 8).Value()   # H - current balance
    $ws.Cells.Item($dst, 12).Value() = $ws.Cells.Item($r, 12).Value()  # L - sale price
    $ws.Cells.Item($dst, 14).Value() = $ws.Cells.Item($r, 14).Value()  # N - number of transactions
}

# --- Step 2: write the new product row 25 (TAMSULIN) ---
$ws.Cells.Item(25, 2).Value()  = "TAMSULIN 0.4MG 28 CAPS"
$ws.Cells.Item(25, 8).Value()  = "2:0"
$ws.Cells.Item(25, 12).Value() = 62
$ws.Cells.Item(25, 14).Value() = "0:2"

# --- Step 3: move the totals row (was 31, now 32) and footer row (was 32, now 33) ---
$oldTotal = $ws.Cells.Item(31, 11).Value()

$footerA = $ws.Cells.Item(32, 1).Value()
$footerF = $ws.Cells.Item(32, 6).Value()
$footerI = $ws.Cells.Item(32, 9).Value()

$ws.Range("K31:N31").UnMerge()
$ws.Range("A32:E32").UnMerge()
$ws.Range("F32:G32").UnMerge()
$ws.Range("I32:N32").UnMerge()

# new row 31 holds the last product (سرنجات 3 سم) - give it its A index and the
# same merges the other product rows have
$ws.Cells.Item(31, 1).Value() = 28
$ws.Range("B31:G31").Merge()
$ws.Range("H31:K31").Merge()
$ws.Range("L31:M31").Merge()

# clear the cells of the old totals/footer rows before relocating their content
$ws.Cells.Item(31, 11).Value() = ""
$ws.Cells.Item(32, 1).Value()  = ""
$ws.Cells.Item(32, 6).Value()  = ""
$ws.Cells.Item(32, 9).Value()  = ""

# totals row moves to 32
$ws.Cells.Item(32, 11).Value() = $oldTotal + 62
$ws.Range("K32:N32").Merge()

# footer row moves to 33
$ws.Cells.Item(33, 1).Value() = $footerA
$ws.Cells.Item(33, 6).Value() = $footerF
$ws.Cells.Item(33, 9).Value() = $footerI
$ws.Range("A33:E33").Merge()
$ws.Range("F33:G33").Merge()
$ws.Range("I33:N33").Merge()

# --- Step 4: row heights ---
$ws.Rows.Item(25).RowHeight = 25.5
$ws.Rows.Item(26).RowHeight = 24.75
$ws.Rows.Item(27).RowHeight = 25.5
$ws.Rows.Item(28).RowHeight = 25.5
$ws.Rows.Item(29).RowHeight = 24.75
$ws.Rows.Item(30).RowHeight = 25.5
$ws.Rows.Item(31).RowHeight = 24.75
$ws.Rows.Item(32).RowHeight = 26.25
$ws.Rows.Item(33).RowHeight = 16.5
